# Updated cryptos list refresh: new Price (col D) and Volume(1h) (col E)
# text values for the coin rows. Percent cells keep their original
# "  +x.xx%  " padding. Price cells that look like a plain decimal number
# (single dot) are forced to stay text via a temporary "@" (Text) number
# format, then the cell style is reset to "Normal" so no stray style index
# is left behind on the cell (matches the source file, where every D/E
# cell is a plain unstyled inline string).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.887.86"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.082.00"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "617.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.363"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.077.97"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.733"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "89.839.63"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "3.648.88"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "3.158.48"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "436.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.15%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.245"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +19.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.175"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.33%  "
$ws.Range("E32").Value = "  +32.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("E34").Value = "  +11.92%  "
$ws.Range("E35").Value = "  +10.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +25.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "483.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("E41").Value = "  -7.19%  "
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.416"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.680"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("E51").Value = "  -0.17%  "
